$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the two "Binance" conversion lines ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.11 = 7695.05 pesos`n✅ 7695.05 pesos = 2.09 = 900.26 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- tasas!N10/O10/N12/O12: updated rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 475
$ws2.Range("O10").Value = 3655.15
$ws2.Range("N12").Value = 3684
$ws2.Range("O12").Value = 431
